$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 1335
$ws.Range("I100").Value = 866.6667
$ws.Range("K100").Value = 866.6667
$ws.Range("M100").Value = -325.6667
$ws.Range("H116").Value = 3210.4
$ws.Range("J116").Value = 3600.5
$ws.Range("L116").Value = 3600.5
$ws.Range("N116").Value = -10484.5
$ws.Range("H132").Value = 8137698
$ws.Range("I132").Value = 9528449
$ws.Range("J132").Value = 24985.334
$ws.Range("K132").Value = 28585347
$ws.Range("L132").Value = 74956.00199999999
$ws.Range("M132").Value = -28582817
$ws.Range("N132").Value = -80016.00199999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1317.7646
$ws.Range("I2").Value = 1055.6364
$ws.Range("J2").Value = 1798.3334
$ws.Range("K2").Value = 1055.6364
$ws.Range("L2").Value = 1798.3334
$ws.Range("M2").Value = -942.6364000000001
$ws.Range("N2").Value = -2024.3334
$ws.Range("H5").Value = 215
$ws.Range("I5").Value = 176.42857
$ws.Range("J5").Value = 350
$ws.Range("K5").Value = 176.42857
$ws.Range("L5").Value = 350
$ws.Range("M5").Value = -64.42857000000001
$ws.Range("N5").Value = -574
$ws.Range("H32").Value = 8805.821
$ws.Range("I32").Value = 6913.7944
$ws.Range("J32").Value = 21362
$ws.Range("K32").Value = 6913.7944
$ws.Range("L32").Value = 21362
$ws.Range("M32").Value = -6626.7944
$ws.Range("N32").Value = -21936
$ws.Range("I61").Value = 52632630
$ws.Range("K61").Value = 52632630
$ws.Range("M61").Value = -52632418
$ws.Range("H63").Value = 18870570
$ws.Range("I63").Value = 2281.5217
$ws.Range("K63").Value = 2281.5217
$ws.Range("M63").Value = -1595.5217
$ws.Range("H66").Value = 18870570
$ws.Range("I66").Value = 2281.5217
$ws.Range("K66").Value = 11407.6085
$ws.Range("M66").Value = -7975.608499999998
$ws.Range("H74").Value = 2510.5925
$ws.Range("I74").Value = 1778.4375
$ws.Range("J74").Value = 3575.5454
$ws.Range("K74").Value = 1778.4375
$ws.Range("L74").Value = 3575.5454
$ws.Range("M74").Value = -904.4375
$ws.Range("N74").Value = -5323.5454
$ws.Range("H77").Value = 2510.5925
$ws.Range("I77").Value = 1778.4375
$ws.Range("J77").Value = 3575.5454
$ws.Range("K77").Value = 8892.1875
$ws.Range("L77").Value = 17877.727
$ws.Range("M77").Value = -4524.1875
$ws.Range("N77").Value = -26613.727
$ws.Range("H110").Value = 2394.7144
$ws.Range("I110").Value = 1487.5
$ws.Range("K110").Value = 1487.5
$ws.Range("M110").Value = 557.5
$ws.Range("H116").Value = 1317.7646
$ws.Range("I116").Value = 1055.6364
$ws.Range("J116").Value = 1798.3334
$ws.Range("K116").Value = 1055.6364
$ws.Range("L116").Value = 1798.3334
$ws.Range("M116").Value = 1238.3636
$ws.Range("N116").Value = -6386.3334
$ws.Range("I136").Value = 52632630
$ws.Range("K136").Value = 157897890
$ws.Range("M136").Value = -157895340

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1317.7646
$ws.Range("I3").Value = 1055.6364
$ws.Range("J3").Value = 1798.3334
$ws.Range("K3").Value = 1055.6364
$ws.Range("L3").Value = 1798.3334
$ws.Range("M3").Value = -941.6364000000001
$ws.Range("N3").Value = -2026.3334
$ws.Range("H4").Value = 215
$ws.Range("I4").Value = 176.42857
$ws.Range("J4").Value = 350
$ws.Range("K4").Value = 176.42857
$ws.Range("L4").Value = 350
$ws.Range("M4").Value = -61.42857000000001
$ws.Range("N4").Value = -580
$ws.Range("H94").Value = 10417042
$ws.Range("I94").Value = 10869945
$ws.Range("K94").Value = 10869945
$ws.Range("M94").Value = -10869494
$ws.Range("H105").Value = 76924490
$ws.Range("I105").Value = 111112230
$ws.Range("J105").Value = 2050
$ws.Range("K105").Value = 111112230
$ws.Range("L105").Value = 2050
$ws.Range("M105").Value = -111110483
$ws.Range("N105").Value = -5544
$ws.Range("H139").Value = 37520
$ws.Range("J139").Value = 37520
$ws.Range("L139").Value = 37520
$ws.Range("N139").Value = -47800

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 142858910
$ws.Range("H113").Value = 142858910
$ws.Range("H141").Value = 452441.44
$ws.Range("J141").Value = 486013.84
$ws.Range("L141").Value = 486013.84
$ws.Range("N141").Value = -496373.84

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 212.66667
$ws.Range("I38").Value = 169
$ws.Range("J38").Value = 300
$ws.Range("K38").Value = 507
$ws.Range("L38").Value = 900
$ws.Range("M38").Value = -160
$ws.Range("N38").Value = -1594
$ws.Range("H98").Value = 1479.4445
$ws.Range("I98").Value = 1926
$ws.Range("J98").Value = 586.3333
$ws.Range("K98").Value = 5778
$ws.Range("L98").Value = 1758.9999
$ws.Range("M98").Value = -4280
$ws.Range("N98").Value = -4754.9999
$ws.Range("H131").Value = 28575172
$ws.Range("I131").Value = 83333750
$ws.Range("J131").Value = 5478.087
$ws.Range("K131").Value = 250001250
$ws.Range("L131").Value = 16434.261
$ws.Range("M131").Value = -249996210
$ws.Range("N131").Value = -26514.261

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2374.279
$ws.Range("I132").Value = 2036
$ws.Range("K132").Value = 6108
$ws.Range("M132").Value = -3578

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1246.5
$ws.Range("J61").Value = 1441.6
$ws.Range("L61").Value = 1441.6
$ws.Range("N61").Value = -1845.6
$ws.Range("H113").Value = 1246.5
$ws.Range("J113").Value = 1441.6
$ws.Range("L113").Value = 1441.6
$ws.Range("N113").Value = -5781.6
$ws.Range("H122").Value = 10871080
$ws.Range("I122").Value = 14707343
$ws.Range("J122").Value = 1667.3334
$ws.Range("K122").Value = 44122029
$ws.Range("L122").Value = 5002.0002
$ws.Range("M122").Value = -44119579
$ws.Range("N122").Value = -9902.0002
$ws.Range("H132").Value = 2956.818
$ws.Range("I132").Value = 2895.6365
$ws.Range("J132").Value = 3018
$ws.Range("K132").Value = 8686.9095
$ws.Range("L132").Value = 9054
$ws.Range("M132").Value = -6156.9095
$ws.Range("N132").Value = -14114

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 499.92856
$ws.Range("I113").Value = 345
$ws.Range("J113").Value = 778.8
$ws.Range("K113").Value = 1035
$ws.Range("L113").Value = 2336.4
$ws.Range("M113").Value = 1135
$ws.Range("N113").Value = -6676.4
$ws.Range("H126").Value = 111112670
$ws.Range("I126").Value = 333333900
$ws.Range("J126").Value = 2051.6667
$ws.Range("K126").Value = 1000001700
$ws.Range("L126").Value = 6155.000100000001
$ws.Range("M126").Value = -999999230
$ws.Range("N126").Value = -11095.0001
$ws.Range("H136").Value = 1243.3077
$ws.Range("I136").Value = 948.2941
$ws.Range("K136").Value = 2844.8823
$ws.Range("M136").Value = -294.8822999999998
